# The document has two blank paragraphs (each consisting only of a
# paragraph mark, formatted bold / rFonts cstheme="minorHAnsi") sitting
# directly after the first table and before the run of similarly blank
# paragraphs that precede the second table. Remove those first two blank
# paragraphs, leaving the rest of the run (and everything else) untouched.

$d = $word.ActiveDocument

# Anchor on the end of the first table - the blank paragraphs start right
# after it.
$t1 = $d.Tables.Item(1)
$pos = $t1.Range.End

# Deleting a single character here removes that paragraph's mark, which
# merges it away while the following paragraph slides into its place.
# Doing this twice (re-reading $pos each time, since position 780 stays
# the same insertion point) removes the first two blank paragraphs.
$rng = $d.Range($pos, $pos + 1)
$rng.Delete()

$rng = $d.Range($pos, $pos + 1)
$rng.Delete()
